{"js": "// Replace each two-digit-by-two-digit multiplication \"problem=answer\"\n// string in the document with its updated value, per the commit diff.\n// Every old string is unique within the document, so a plain body-wide\n// search + full-text replace per pair is safe and keeps the original\n// run formatting (rFonts/sz) untouched.\nconst replacements = [\n  [\"29\u00d714=406\", \"75\u00d777=5775\"],\n  [\"38\u00d744=1672\", \"68\u00d735=2380\"],\n  [\"21\u00d746=966\", \"69\u00d726=1794\"],\n  [\"96\u00d785=8160\", \"61\u00d753=3233\"],\n  [\"75\u00d765=4875\", \"32\u00d757=1824\"],\n  [\"87\u00d751=4437\", \"90\u00d746=4140\"],\n  [\"68\u00d715=1020\", \"34\u00d762=2108\"],\n  [\"38\u00d780=3040\", \"45\u00d718=810\"],\n  [\"43\u00d763=2709\", \"95\u00d740=3800\"],\n  [\"49\u00d719=931\", \"78\u00d726=2028\"],\n  [\"51\u00d771=3621\", \"72\u00d756=4032\"],\n  [\"67\u00d763=4221\", \"51\u00d718=918\"],\n  [\"81\u00d777=6237\", \"54\u00d775=4050\"],\n  [\"20\u00d752=1040\", \"38\u00d749=1862\"],\n  [\"38\u00d713=494\", \"75\u00d719=1425\"],\n  [\"22\u00d789=1958\", \"95\u00d716=1520\"],\n  [\"94\u00d730=2820\", \"83\u00d766=5478\"],\n  [\"74\u00d762=4588\", \"79\u00d779=6241\"],\n  [\"27\u00d754=1458\", \"38\u00d745=1710\"],\n  [\"26\u00d770=1820\", \"97\u00d748=4656\"],\n  [\"36\u00d777=2772\", \"68\u00d721=1428\"],\n  [\"84\u00d732=2688\", \"88\u00d793=8184\"],\n  [\"87\u00d773=6351\", \"94\u00d711=1034\"],\n  [\"61\u00d736=2196\", \"88\u00d717=1496\"],\n  [\"29\u00d770=2030\", \"38\u00d777=2926\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-by-two-digit multiplication \"problem=answer\"\n# string in the document with its updated value, per the commit diff.\n# Every old string is unique within the document, so Find/Replace across\n# the whole document body (wdReplaceAll) for each pair is safe and\n# leaves the surrounding run formatting (rFonts/sz) untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"29\u00d714=406\", \"75\u00d777=5775\"),\n  @(\"38\u00d744=1672\", \"68\u00d735=2380\"),\n  @(\"21\u00d746=966\", \"69\u00d726=1794\"),\n  @(\"96\u00d785=8160\", \"61\u00d753=3233\"),\n  @(\"75\u00d765=4875\", \"32\u00d757=1824\"),\n  @(\"87\u00d751=4437\", \"90\u00d746=4140\"),\n  @(\"68\u00d715=1020\", \"34\u00d762=2108\"),\n  @(\"38\u00d780=3040\", \"45\u00d718=810\"),\n  @(\"43\u00d763=2709\", \"95\u00d740=3800\"),\n  @(\"49\u00d719=931\", \"78\u00d726=2028\"),\n  @(\"51\u00d771=3621\", \"72\u00d756=4032\"),\n  @(\"67\u00d763=4221\", \"51\u00d718=918\"),\n  @(\"81\u00d777=6237\", \"54\u00d775=4050\"),\n  @(\"20\u00d752=1040\", \"38\u00d749=1862\"),\n  @(\"38\u00d713=494\", \"75\u00d719=1425\"),\n  @(\"22\u00d789=1958\", \"95\u00d716=1520\"),\n  @(\"94\u00d730=2820\", \"83\u00d766=5478\"),\n  @(\"74\u00d762=4588\", \"79\u00d779=6241\"),\n  @(\"27\u00d754=1458\", \"38\u00d745=1710\"),\n  @(\"26\u00d770=1820\", \"97\u00d748=4656\"),\n  @(\"36\u00d777=2772\", \"68\u00d721=1428\"),\n  @(\"84\u00d732=2688\", \"88\u00d793=8184\"),\n  @(\"87\u00d773=6351\", \"94\u00d711=1034\"),\n  @(\"61\u00d736=2196\", \"88\u00d717=1496\"),\n  @(\"29\u00d770=2030\", \"38\u00d777=2926\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $rng = $d.Content\n  $find = $rng.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute(\n    $oldText,   # FindText\n    $true,      # MatchCase\n    $false,     # MatchWholeWord\n    $false,     # MatchWildcards\n    $false,     # MatchSoundsLike\n    $false,     # MatchAllWordForms\n    $true,      # Forward\n    1,          # Wrap (wdFindContinue)\n    $false,     # Format\n    $newText,   # ReplaceWith\n    2           # Replace (wdReplaceAll)\n  ) | Out-Null\n}\n"}
